$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.09179766666666667
$ws.Range("H2").Value = 0.275393
$ws.Range("M2").Value = 1.135923333333333
$ws.Range("N2").Value = 3.40777
$ws.Range("O2").Value = 0.06998805686568385
$ws.Range("P2").Value = 0.06998805686568385
$ws.Range("Q2").Value = 0.1042751115122222
$ws.Range("R2").Value = 0.9384760036099999
$ws.Range("S2").Value = 0.06998805686568385
$ws.Range("T2").Value = 0.06998805686568385

# Row 3
$ws.Range("G3").Value = 0.09179766666666667
$ws.Range("H3").Value = 0.275393
$ws.Range("O3").Value = 0.4371656037403091
$ws.Range("P3").Value = 0.437165603740309
$ws.Range("Q3").Value = 0.6513324432883333
$ws.Range("R3").Value = 5.861991989594999
$ws.Range("S3").Value = 0.4371656037403091
$ws.Range("T3").Value = 0.437165603740309

# Row 4
$ws.Range("G4").Value = 0.09179766666666667
$ws.Range("H4").Value = 0.275393
$ws.Range("M4").Value = 7.999016999999999
$ws.Range("N4").Value = 23.997051
$ws.Range("O4").Value = 0.4928463393940071
$ws.Range("P4").Value = 0.4928463393940071
$ws.Range("Q4").Value = 0.734291096227
$ws.Range("R4").Value = 6.608619866043
$ws.Range("S4").Value = 0.4928463393940071
$ws.Range("T4").Value = 0.4928463393940071
